# Append the 2025-02-24 22:11 resale-number update as new row 87.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

# --- Text columns (A-D) ---------------------------------------------------
# Force text storage so date/time-looking strings and the zero-padded
# "08" week number aren't auto-converted to a date serial / number by Excel.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-02-24"
$ws.Cells.Item($row, 2).Value = "22:11:36"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "08"

# Drop back to the default "Normal" style so these new cells don't pick up
# an explicit number-format style index (matches the unstyled cells used
# throughout the rest of the data rows).
$textRange.Style = "Normal"

# --- Numeric columns (E-T) --------------------------------------------------
$ws.Cells.Item($row, 5).Value  = 130720
$ws.Cells.Item($row, 6).Value  = 141927
$ws.Cells.Item($row, 7).Value  = 172841
$ws.Cells.Item($row, 8).Value  = 158991
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 146625
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193851
$ws.Cells.Item($row, 14).Value = 115476
$ws.Cells.Item($row, 15).Value = 46528
$ws.Cells.Item($row, 16).Value = 29420
$ws.Cells.Item($row, 17).Value = 68856
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48926
$ws.Cells.Item($row, 20).Value = -1
